$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K11").Value = 0.2305062539156956
$ws.Range("J12").Value = 0.2193215401759246
$ws.Range("I13").Value = 0.2109873117084238
$ws.Range("H14").Value = 0.2077622620068982
$ws.Range("G15").Value = 0.1868984584576193
$ws.Range("F16").Value = 0.2101374940836094
$ws.Range("E17").Value = 0.2201756597651073
$ws.Range("D18").Value = 0.1085991175498651
$ws.Range("C19").Value = 0.130019622424466
$ws.Range("B20").Value = 0.3662627537369125
